# Commit message: "Change script to manager, Update Process book, Put down a to-do list"
#
# The document has 5 consecutive blank paragraphs right after the
# "As we analyze further..." paragraph and right before the
# "Exploratory Data Analysis: ..." heading. This change fills those 5
# blank paragraphs in with reflection / to-do text, matching the
# formatting already defined by each paragraph mark (Arial, 10.5pt /
# half-point size 21, single underline, black, no extra character
# spacing or position offset).

$d = $word.ActiveDocument

function Set-ParaText($paragraph, [string]$text) {
    $r = $paragraph.Range
    $r.Text = $text
    $r.Font.Name = "Arial"
    $r.Font.NameFarEast = "Arial"
    $r.Font.NameOther = "Arial"
    $r.Font.NameBi = "Arial"
    $r.Font.Size = 10.5
    $r.Font.Underline = 1
    $r.Font.Color = 0
    $r.Font.Spacing = 0
    $r.Font.Position = 0
}

# Locate the anchor paragraph ("As we analyze further...") via Find so we
# are not dependent on a brittle hard-coded paragraph number.
$findRange = $d.Content
$ok = $findRange.Find.Execute("As we analyze further", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    Write-Host "ERROR: anchor paragraph not found"
}

$anchor = $findRange.Paragraphs.Item(1)
$anchor = $d.Paragraphs.Item($anchor.Index)

$blank1 = $anchor.Next()
$blank2 = $blank1.Next()
$blank3 = $blank2.Next()
$blank4 = $blank3.Next()
$blank5 = $blank4.Next()

Set-ParaText $blank1 "A scaling factor is calculated from visible tree and overall tree to transform the treenode. We can get rid of the lagging behavior of the visualization in this case. "
Set-ParaText $blank2 "We separate different classes of the visualizations so that these visualizations do not know and depend on each other. There is a manager script that knows about the behavior of all the classes and all the interactions are done in this class. "
Set-ParaText $blank3 "We add some animation when updating the tree structure. "
Set-ParaText $blank4 "There is some lagging behavior when first constructing the tree. We are not sure right whether that can be got rid of. "
Set-ParaText $blank5 "Several things that may be done: 1. Changing the size of the tooltip. 2. Labeling the axis of the plot.  3. Interaction for user's input (Slider would be better than buttons)"
